# Add new worksheet 'Vert' to cater vertical data types (orientation parameter)
$wb = $excel.ActiveWorkbook

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Vert"

$newSheet.Range("A1").Value = "name"
$newSheet.Range("B1").Value = "value1"
$newSheet.Range("A2").Value = "desc"
$newSheet.Range("B2").Value = "value2"
$newSheet.Range("A3").Value = "cidr_block"
$newSheet.Range("B3").Value = "value3"

$newSheet.Activate()
